$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 2245
$ws.Cells.Item(7, 6).Value = 8161
$ws.Cells.Item(8, 6).Value = 105
$ws.Cells.Item(9, 6).Value = 9
$ws.Cells.Item(10, 6).Value = 1799
$ws.Cells.Item(11, 6).Value = 1595
$ws.Cells.Item(12, 6).Value = 1324
$ws.Cells.Item(13, 6).Value = 209
$ws.Cells.Item(14, 6).Value = 4389
$ws.Cells.Item(15, 6).Value = 6164
$ws.Cells.Item(16, 6).Value = 761
$ws.Cells.Item(17, 6).Value = 57
$ws.Cells.Item(18, 6).Value = 1198
$ws.Cells.Item(20, 6).Value = 469
$ws.Cells.Item(21, 6).Value = 6427
$ws.Cells.Item(24, 6).Value = 4359
$ws.Cells.Item(25, 6).Value = 301
$ws.Cells.Item(26, 6).Value = 717
$ws.Cells.Item(27, 6).Value = 2015
$ws.Cells.Item(29, 6).Value = 339
$ws.Cells.Item(31, 6).Value = 35
$ws.Cells.Item(32, 6).Value = 48
$ws.Cells.Item(34, 6).Value = 82
$ws.Cells.Item(35, 6).Value = 333
$ws.Cells.Item(36, 6).Value = 1187
$ws.Cells.Item(37, 6).Value = 1901
$ws.Cells.Item(38, 6).Value = 134
$ws.Cells.Item(39, 6).Value = 440
$ws.Cells.Item(41, 6).Value = 1204
$ws.Cells.Item(43, 6).Value = 70
$ws.Cells.Item(44, 6).Value = 1147
$ws.Cells.Item(46, 6).Value = 82
$ws.Cells.Item(47, 6).Value = 189
$ws.Cells.Item(48, 6).Value = 28

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 22
$ws.Cells.Item(11, 6).Value = 414
$ws.Cells.Item(14, 6).Value = 115
$ws.Cells.Item(18, 6).Value = 191
$ws.Cells.Item(29, 6).Value = 129
$ws.Cells.Item(37, 6).Value = 17

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 1591
$ws.Cells.Item(7, 6).Value = 491
$ws.Cells.Item(9, 6).Value = 1054
$ws.Cells.Item(10, 6).Value = 1144
$ws.Cells.Item(11, 6).Value = 1447
$ws.Cells.Item(12, 6).Value = 1811
$ws.Cells.Item(13, 6).Value = 300
$ws.Cells.Item(14, 6).Value = 151

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1591
$ws.Cells.Item(6, 6).Value = 491
$ws.Cells.Item(8, 6).Value = 2245
$ws.Cells.Item(9, 6).Value = 1054
$ws.Cells.Item(10, 6).Value = 9
$ws.Cells.Item(11, 6).Value = 22
$ws.Cells.Item(12, 6).Value = 1595
$ws.Cells.Item(13, 6).Value = 1447
$ws.Cells.Item(14, 6).Value = 1324
$ws.Cells.Item(16, 6).Value = 209
$ws.Cells.Item(17, 6).Value = 1811
$ws.Cells.Item(18, 6).Value = 4389
$ws.Cells.Item(19, 6).Value = 300
$ws.Cells.Item(21, 6).Value = 414
$ws.Cells.Item(22, 6).Value = 761
$ws.Cells.Item(23, 6).Value = 57
$ws.Cells.Item(24, 6).Value = 1199
$ws.Cells.Item(26, 6).Value = 469
$ws.Cells.Item(27, 6).Value = 6427
$ws.Cells.Item(29, 6).Value = 151
$ws.Cells.Item(31, 6).Value = 4359
$ws.Cells.Item(32, 6).Value = 301
$ws.Cells.Item(33, 6).Value = 2015
$ws.Cells.Item(35, 6).Value = 339
$ws.Cells.Item(38, 6).Value = 191
$ws.Cells.Item(39, 6).Value = 82
$ws.Cells.Item(40, 6).Value = 333
$ws.Cells.Item(41, 6).Value = 1901
$ws.Cells.Item(42, 6).Value = 134
$ws.Cells.Item(43, 6).Value = 440
$ws.Cells.Item(44, 6).Value = 1204
$ws.Cells.Item(48, 6).Value = 1147
$ws.Cells.Item(49, 6).Value = 189
